# continuing refactor of power calculation
#
# Triggers sheet: the two old duplicate "Compression" columns (G & H, all
# zeros) are collapsed into a single column that is moved to the end, and
# the "Link Efficiency (J/bit)" / "Op Efficiency (J/op)" columns shift left
# to take their place (G & H). Op Efficiency values are also updated for the
# Intermediate->Global and Global->Output rows.

$wb = $excel.ActiveWorkbook
$detectors = $wb.Worksheets.Item("Detectors")
$triggers  = $wb.Worksheets.Item("Triggers")

# ---------------------------------------------------------------------
# Triggers sheet data rewrite
# ---------------------------------------------------------------------

# Header row (row 1): new column order is
# Name | Output | Data (bytes) | Reduction | Skill mean | Skill variance |
# Link Efficiency (J/bit) | Op Efficiency (J/op) | Compression
$triggers.Cells.Item(1, 7).Value = "Link Efficiency (J/bit)"
$triggers.Cells.Item(1, 8).Value = "Op Efficiency (J/op)"
$triggers.Cells.Item(1, 9).Value = "Compression"

# Data rows: shift Link Efficiency (old col I) into col G, move the Op
# Efficiency values (old col J) into col H (with the two updated values),
# and the Compression value (old col G, always 0) into col I.
$linkEfficiency = 0.000000000025

$triggers.Cells.Item(2, 7).Value = $linkEfficiency
$triggers.Cells.Item(2, 8).Value = 0
$triggers.Cells.Item(2, 9).Value = 0

$triggers.Cells.Item(3, 7).Value = $linkEfficiency
$triggers.Cells.Item(3, 8).Value = 0
$triggers.Cells.Item(3, 9).Value = 0

$triggers.Cells.Item(4, 7).Value = $linkEfficiency
$triggers.Cells.Item(4, 8).Value = 0
$triggers.Cells.Item(4, 9).Value = 0

$triggers.Cells.Item(5, 7).Value = $linkEfficiency
$triggers.Cells.Item(5, 8).Value = 0
$triggers.Cells.Item(5, 9).Value = 0

$triggers.Cells.Item(6, 7).Value = $linkEfficiency
$triggers.Cells.Item(6, 8).Value = 0.003
$triggers.Cells.Item(6, 9).Value = 0

$triggers.Cells.Item(7, 7).Value = $linkEfficiency
$triggers.Cells.Item(7, 8).Value = 16
$triggers.Cells.Item(7, 9).Value = 0

$triggers.Cells.Item(8, 7).Value = $linkEfficiency
$triggers.Cells.Item(8, 8).Value = 0
$triggers.Cells.Item(8, 9).Value = 0

# The old column J (Op Efficiency) is no longer used - drop it so the used
# range shrinks back from A1:J8 to A1:I8.
$triggers.Columns.Item(10).Clear()

# ---------------------------------------------------------------------
# Styling
# ---------------------------------------------------------------------

# Whole table uses the bold/explicit-black font (fontId 1) now, matching
# the header row's existing look.
$triggers.Range("A1:I8").Font.Color = 0

# The Link Efficiency column (now G) additionally keeps its scientific
# number format on the data rows.
$triggers.Range("G2:G8").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------

$detectors.Range("G1").Select()
$triggers.Range("F23").Select()
